# RPA datasets push 2024-06-08
# Three IPO listings that previously had no confirmed offering price ("-")
# now have their 확정공모가 (confirmed offering price) filled in:
#   row 15 (에이치엠씨아이비스팩7호)   -> 2000
#   row 17 (미래에셋비전스팩5호)       -> 2000
#   row 21 (한중엔시에스)              -> 30000
# These values must be written as text (matching the existing D-column
# cells, which are all shared-string "-" placeholders / "2000" text,
# not numbers) - so we use Excel's leading-apostrophe text prefix, then
# clear the resulting cell style back to Normal so no stray formatting
# is left behind on the cell itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Value = "'2000"
$ws.Range("D17").Value = "'2000"
$ws.Range("D21").Value = "'30000"

$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D21").Style = "Normal"
